$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PIR")
$rows = @(
    @("2026-01-28", "14:54:29", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "14:54:33", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "14:54:38", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "14:54:43", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "14:54:48", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "14:54:53", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "14:54:58", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "14:55:03", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "14:55:08", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "14:55:13", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "14:55:18", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "14:55:23", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "14:55:28", "14:00", "Bathroom", "No Motion", "Inactive")
)
$startRow = 84
$r = $startRow
foreach ($rowData in $rows) {
    $ws.Cells.Item($r, 1).Value = "'" + $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    if ($rowData[4] -like '*%') {
        $ws.Cells.Item($r, 5).Value = "'" + $rowData[4]
    } else {
        $ws.Cells.Item($r, 5).Value = $rowData[4]
    }
    $ws.Cells.Item($r, 6).Value = $rowData[5]
    $r = $r + 1
}

$ws = $wb.Worksheets.Item("Humidity")
$rows = @(
    @("2026-01-28", "14:54:28", "14:00", "Bathroom", "88.3%", "Active"),
    @("2026-01-28", "14:54:40", "14:00", "Bathroom", "87.4%", "Active"),
    @("2026-01-28", "14:54:44", "14:00", "Bathroom", "88.3%", "Active"),
    @("2026-01-28", "14:54:48", "14:00", "Bathroom", "87.4%", "Active"),
    @("2026-01-28", "14:54:52", "14:00", "Bathroom", "88.3%", "Active"),
    @("2026-01-28", "14:55:00", "14:00", "Bathroom", "87.4%", "Active"),
    @("2026-01-28", "14:55:04", "14:00", "Bathroom", "88.3%", "Active"),
    @("2026-01-28", "14:55:08", "14:00", "Bathroom", "87.4%", "Active"),
    @("2026-01-28", "14:55:13", "14:00", "Bathroom", "88.3%", "Active"),
    @("2026-01-28", "14:55:16", "14:00", "Bathroom", "88.3%", "Active"),
    @("2026-01-28", "14:55:20", "14:00", "Bathroom", "88.3%", "Active"),
    @("2026-01-28", "14:55:25", "14:00", "Bathroom", "88.3%", "Active")
)
$startRow = 80
$r = $startRow
foreach ($rowData in $rows) {
    $ws.Cells.Item($r, 1).Value = "'" + $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    if ($rowData[4] -like '*%') {
        $ws.Cells.Item($r, 5).Value = "'" + $rowData[4]
    } else {
        $ws.Cells.Item($r, 5).Value = $rowData[4]
    }
    $ws.Cells.Item($r, 6).Value = $rowData[5]
    $r = $r + 1
}

$ws = $wb.Worksheets.Item("Temperature")
$rows = @(
    @("2026-01-28", "14:54:29", "14:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "14:54:40", "14:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "14:54:44", "14:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "14:54:48", "14:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "14:54:52", "14:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "14:55:01", "14:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "14:55:05", "14:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "14:55:09", "14:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "14:55:13", "14:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "14:55:17", "14:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "14:55:21", "14:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "14:55:25", "14:00", "Bathroom", "22.8C", "Active")
)
$startRow = 80
$r = $startRow
foreach ($rowData in $rows) {
    $ws.Cells.Item($r, 1).Value = "'" + $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    if ($rowData[4] -like '*%') {
        $ws.Cells.Item($r, 5).Value = "'" + $rowData[4]
    } else {
        $ws.Cells.Item($r, 5).Value = $rowData[4]
    }
    $ws.Cells.Item($r, 6).Value = $rowData[5]
    $r = $r + 1
}

$ws = $wb.Worksheets.Item("Proximity")
$rows = @(
    @("2026-01-28", "14:55:12", "14:00", "Bathroom Door", "ENTER", "User ENTERED Bathroom"),
    @("2026-01-28", "14:55:12", "14:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-01-28", "14:55:15", "14:00", "Living Room Main Door", "EXIT", "User EXITED Living Room Main Door"),
    @("2026-01-28", "14:55:24", "14:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-01-28", "14:55:24", "14:00", "Bathroom Door", "EXIT", "User EXITED Bathroom")
)
$startRow = 19
$r = $startRow
foreach ($rowData in $rows) {
    $ws.Cells.Item($r, 1).Value = "'" + $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    if ($rowData[4] -like '*%') {
        $ws.Cells.Item($r, 5).Value = "'" + $rowData[4]
    } else {
        $ws.Cells.Item($r, 5).Value = $rowData[4]
    }
    $ws.Cells.Item($r, 6).Value = $rowData[5]
    $r = $r + 1
}

$ws = $wb.Worksheets.Item("mmWave")
$rows = @(
    @("2026-01-28", "14:54:31", "14:00", "Living Room", "Presence Detected", "Active")
)
$startRow = 4
$r = $startRow
foreach ($rowData in $rows) {
    $ws.Cells.Item($r, 1).Value = "'" + $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    if ($rowData[4] -like '*%') {
        $ws.Cells.Item($r, 5).Value = "'" + $rowData[4]
    } else {
        $ws.Cells.Item($r, 5).Value = $rowData[4]
    }
    $ws.Cells.Item($r, 6).Value = $rowData[5]
    $r = $r + 1
}

$ws = $wb.Worksheets.Item("Camera")
$rows = @(
    @("2026-01-28", "14:55:14", "14:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-01-28", "14:55:25", "14:00", "Living Room Main Door", "Image Captured", "Active")
)
$startRow = 10
$r = $startRow
foreach ($rowData in $rows) {
    $ws.Cells.Item($r, 1).Value = "'" + $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    if ($rowData[4] -like '*%') {
        $ws.Cells.Item($r, 5).Value = "'" + $rowData[4]
    } else {
        $ws.Cells.Item($r, 5).Value = $rowData[4]
    }
    $ws.Cells.Item($r, 6).Value = $rowData[5]
    $r = $r + 1
}

